$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = -7.333
$ws.Range("C9").Value = -10.925
$ws.Range("D12").Value = -7.697999999999999
$ws.Range("C13").Value = -13.305
$ws.Range("D14").Value = -7.959999999999999
$ws.Range("C16").Value = -13.095
$ws.Range("C18").Value = -11.305
$ws.Range("D19").Value = -8.1
$ws.Range("C20").Value = -12.578
$ws.Range("C26").Value = -12.476
$ws.Range("D26").Value = -8.126999999999999
$ws.Range("C27").Value = -12.84
$ws.Range("D27").Value = -8.433
$ws.Range("C29").Value = -12.34
$ws.Range("D29").Value = -7.419
$ws.Range("C35").Value = -12.377
$ws.Range("C36").Value = -12.428
$ws.Range("D37").Value = -8.06
$ws.Range("D38").Value = -7.722
$ws.Range("C45").Value = -13.009
$ws.Range("D47").Value = -7.502
$ws.Range("D51").Value = -8.311000000000002
$ws.Range("D52").Value = -7.601000000000001
$ws.Range("C55").Value = -13.44
$ws.Range("D55").Value = -8.409
$ws.Range("C57").Value = -13.563
$ws.Range("C69").Value = -11.269
$ws.Range("D69").Value = -7.243
$ws.Range("D70").Value = -7.323000000000002
$ws.Range("C76").Value = -13.42
$ws.Range("D76").Value = -7.276999999999999
$ws.Range("C78").Value = -12.52
$ws.Range("D81").Value = -7.843999999999999
$ws.Range("C82").Value = -11.888
$ws.Range("C83").Value = -13.184
$ws.Range("D83").Value = -8.441999999999998
$ws.Range("C93").Value = -11.724
$ws.Range("D94").Value = -7.455999999999999
$ws.Range("C97").Value = -12.028
$ws.Range("D100").Value = -8.293999999999999
$ws.Range("D102").Value = -7.833
